$d = $word.ActiveDocument

# Select the whole document (the only content is a single, empty
# paragraph whose mark carries direct character formatting: sz/szCs
# 24 half-points i.e. 12pt) and clear that direct formatting, the same
# way a user would by selecting everything and pressing "Clear
# Formatting". This removes the <w:rPr> (sz/szCs) held on the
# paragraph mark's <w:pPr>, leaving a plain empty paragraph.
$sel = $word.Selection
$sel.WholeStory()
$sel.ClearFormatting()
